$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = [double]"146.620486"
$ws.Range("H2").Value = [double]"439.861458"
$ws.Range("I2").Value = [double]"0.3983053592962091"
$ws.Range("J2").Value = [double]"0.3983053592962091"
$ws.Range("O2").Value = [double]"0.02773017886769741"
$ws.Range("P2").Value = [double]"0.02773017886769741"
$ws.Range("Q2").Value = [double]"7.615663536821333"
$ws.Range("R2").Value = [double]"68.54097183139199"
$ws.Range("S2").Value = [double]"0.01104507885724636"
$ws.Range("T2").Value = [double]"0.01104507885724636"

# Row 3
$ws.Range("G3").Value = [double]"146.620486"
$ws.Range("H3").Value = [double]"439.861458"
$ws.Range("I3").Value = [double]"0.3983053592962091"
$ws.Range("J3").Value = [double]"0.3983053592962091"
$ws.Range("M3").Value = [double]"1.821156333333333"
$ws.Range("N3").Value = [double]"5.463469"
$ws.Range("O3").Value = [double]"0.9722698211323025"
$ws.Range("P3").Value = [double]"0.9722698211323026"
$ws.Range("Q3").Value = [double]"267.0188266753113"
$ws.Range("R3").Value = [double]"2403.169440077802"
$ws.Range("S3").Value = [double]"0.3872602804389627"
$ws.Range("T3").Value = [double]"0.3872602804389628"

# Row 4
$ws.Range("I4").Value = [double]"0.534552907532962"
$ws.Range("J4").Value = [double]"0.5345529075329621"
$ws.Range("O4").Value = [double]"0.02773017886769741"
$ws.Range("P4").Value = [double]"0.02773017886769741"
$ws.Range("S4").Value = [double]"0.01482324774013675"
$ws.Range("T4").Value = [double]"0.01482324774013676"

# Row 5
$ws.Range("I5").Value = [double]"0.534552907532962"
$ws.Range("J5").Value = [double]"0.5345529075329621"
$ws.Range("M5").Value = [double]"1.821156333333333"
$ws.Range("N5").Value = [double]"5.463469"
$ws.Range("O5").Value = [double]"0.9722698211323025"
$ws.Range("P5").Value = [double]"0.9722698211323026"
$ws.Range("Q5").Value = [double]"358.3574431876499"
$ws.Range("R5").Value = [double]"3225.216988688849"
$ws.Range("S5").Value = [double]"0.5197296597928253"
$ws.Range("T5").Value = [double]"0.5197296597928254"

# Row 6
$ws.Range("G6").Value = [double]"24.174389"
$ws.Range("H6").Value = [double]"72.523167"
$ws.Range("I6").Value = [double]"0.0656715098899026"
$ws.Range("J6").Value = [double]"0.0656715098899026"
$ws.Range("O6").Value = [double]"0.02773017886769741"
$ws.Range("P6").Value = [double]"0.02773017886769741"
$ws.Range("Q6").Value = [double]"1.255649997178667"
$ws.Range("R6").Value = [double]"11.300849974608"
$ws.Range("S6").Value = [double]"0.001821082715758759"
$ws.Range("T6").Value = [double]"0.001821082715758759"

# Row 7
$ws.Range("G7").Value = [double]"24.174389"
$ws.Range("H7").Value = [double]"72.523167"
$ws.Range("I7").Value = [double]"0.0656715098899026"
$ws.Range("J7").Value = [double]"0.0656715098899026"
$ws.Range("M7").Value = [double]"1.821156333333333"
$ws.Range("N7").Value = [double]"5.463469"
$ws.Range("O7").Value = [double]"0.9722698211323025"
$ws.Range("P7").Value = [double]"0.9722698211323026"
$ws.Range("Q7").Value = [double]"44.02534163181367"
$ws.Range("R7").Value = [double]"396.228074686323"
$ws.Range("S7").Value = [double]"0.06385042717414384"
$ws.Range("T7").Value = [double]"0.06385042717414384"

# Row 8
$ws.Range("E8").Value = [double]"3"
$ws.Range("F8").Value = [double]"1"
$ws.Range("G8").Value = [double]"0.541205"
$ws.Range("H8").Value = [double]"1.623615"
$ws.Range("I8").Value = [double]"0.001470223280926138"
$ws.Range("J8").Value = [double]"0.001470223280926138"
$ws.Range("O8").Value = [double]"0.02773017886769741"
$ws.Range("P8").Value = [double]"0.02773017886769741"
$ws.Range("Q8").Value = [double]"0.02811090930666667"
$ws.Range("R8").Value = [double]"0.25299818376"
$ws.Range("S8").Value = [double]"4.076955455553474E-05"
$ws.Range("T8").Value = [double]"4.076955455553474E-05"

# Row 9
$ws.Range("E9").Value = [double]"3"
$ws.Range("F9").Value = [double]"1"
$ws.Range("G9").Value = [double]"0.541205"
$ws.Range("H9").Value = [double]"1.623615"
$ws.Range("I9").Value = [double]"0.001470223280926138"
$ws.Range("J9").Value = [double]"0.001470223280926138"
$ws.Range("M9").Value = [double]"1.821156333333333"
$ws.Range("N9").Value = [double]"5.463469"
$ws.Range("O9").Value = [double]"0.9722698211323025"
$ws.Range("P9").Value = [double]"0.9722698211323026"
$ws.Range("Q9").Value = [double]"0.9856189133816667"
$ws.Range("R9").Value = [double]"8.870570220435001"
$ws.Range("S9").Value = [double]"0.001429453726370603"
$ws.Range("T9").Value = [double]"0.001429453726370603"
